$wb = $excel.ActiveWorkbook

$oldSheetName = "浙江杭州三墩地铁站店_原始凭证"
$newSheetName = "浙江杭州西湖三墩地铁站店_原始凭证"
$newStoreName = "浙江杭州西湖三墩地铁站店"
$otherSheetName = "浙江杭州滨江中南乐游城店_原始凭证"

# Rename the worksheet.
$ws = $wb.Worksheets.Item($oldSheetName)
$ws.Name = $newSheetName

# Update the store-name cell on that sheet (B11) to match the new name.
$ws.Range("B11").Value = $newStoreName

# Renaming the sheet strips the sheet qualifier from any defined name whose
# RefersTo is a broken reference (#REF!), because the name-recalc can't
# re-resolve an already-broken reference to the (new) sheet. Restore the
# sheet-qualified #REF! for both the workbook-global and the worksheet-local
# "当前费率" defined name, same as every other defined name pair in this
# workbook (each has a localSheetId="0" copy pointing at 浙江杭州滨江中南乐游城店_原始凭证
# and a global copy pointing at the renamed sheet).
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*#REF!*") {
        if ($n.Parent.Name -eq $wb.Name) {
            $n.RefersTo = "=" + $newSheetName + "!#REF!"
        } else {
            $n.RefersTo = "=" + $otherSheetName + "!#REF!"
        }
    }
}
